$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-19 12:42:38"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
